# Update the "Förändrad" (C) date for every data row (2-16) from 46070 to 46072
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

2..16 | ForEach-Object {
    $ws.Cells.Item($_, 3).Value = 46072
}

# Rows 9 & 10 swap their data (Beteckning / Datum / Area)
$ws.Range("A9").Value = "A 50997-2025"
$ws.Range("B9").Value = 45946
$ws.Range("G9").Value = 1.5

$ws.Range("A10").Value = "A 5792-2024"
$ws.Range("B10").Value = 45335
$ws.Range("G10").Value = 5.6

# Rows 11 & 13 swap their data
$ws.Range("A11").Value = "A 13651-2023"
$ws.Range("B11").Value = 45006
$ws.Range("G11").Value = 2.2

$ws.Range("A13").Value = "A 7827-2026"
$ws.Range("B13").Value = 46062.63958333333
$ws.Range("G13").Value = 2.1

# Rows 12 & 14 swap their data
$ws.Range("A12").Value = "A 35642-2023"
$ws.Range("B12").Value = 45147
$ws.Range("G12").Value = 1.2

$ws.Range("A14").Value = "A 7814-2026"
$ws.Range("B14").Value = 46062.61388888889
$ws.Range("G14").Value = 1.1
